$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (row 1) for the summer-reporting header changes.
$ws.Range("E1").Value = "Legal Given Name"
$ws.Range("F1").Value = "Birthdate"
$ws.Range("H1").Value = "Ministry Course Code and Level"
$ws.Range("J1").Value = "Final Percent"
$ws.Range("L1").Value = "Credits"

# Header row is now taller to fit the wrapped/longer header text.
$ws.Rows.Item(1).RowHeight = 64

# Selection moved to the H1:L1 block.
$ws.Range("H1:L1").Select()
